# Commit: "Added Hands On Demos - Day 4"
#
# This edit:
#  1. Removes the trailing "Add Section Header in Titlecase" slide
#     (the last slide in the deck) that was left over from the section
#     template.
#  2. On the "Install Database" slide, drops the MariaDB call-out box
#     and re-centers the remaining MySQL call-out box in the middle of
#     the slide (including center-aligning its two leading blank
#     paragraphs).

$p = $ppt.ActivePresentation

# --- 1. Delete the trailing section-header slide --------------------------
$p.Slides.Item($p.Slides.Count).Delete()

# --- 2. Update the "Install Database" slide --------------------------------
$s = $p.Slides.Item(8)

# The MariaDB call-out box is no longer needed -> remove it entirely.
$mariaDb = $s.Shapes.Item("object 4")
$mariaDb.Delete()

# Re-position the remaining MySQL call-out box to the centre of the slide
# and center-align its two leading (blank) paragraphs.
$mysql = $s.Shapes.Item("object 3")
$mysql.Left = 417.9993750787402
$mysql.Top = 251.9832233464567

$tr = $mysql.TextFrame.TextRange
$tr.Paragraphs(1, 1).ParagraphFormat.Alignment = 2
$tr.Paragraphs(2, 1).ParagraphFormat.Alignment = 2
